$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Marshalling paragraph: three small wording tweaks.
# ---------------------------------------------------------------------------

# 1a. "by calling appropriate getMessageSize" -> "by calling getMessageSize"
#     (keep the edit inside the existing "...calling appropriate " run so the
#     spell/grammar proofErr markers bracketing the next "getMessageSize("
#     run stay correctly paired)
$d.Content.Find.Execute(
    "calling appropriate ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "calling ", 2) | Out-Null

# 1b. "Based on the result, a character array is defined" ->
#     "Based on the result, a buffer of character array is defined"
$d.Content.Find.Execute(
    "Based on the result, a character array is defined",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Based on the result, a buffer of character array is defined", 2) | Out-Null

# 1c. "also copied into the character array and then get sent" ->
#     "also copied into the buffer and then get sent"
$d.Content.Find.Execute(
    "also copied into the character array and then get sent",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "also copied into the buffer and then get sent", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert two brand-new paragraphs right after the marshalling paragraph
#    (which ends in "...and then get sent through socket."), describing the
#    unmarshalling process, before the existing blank separator paragraph.
# ---------------------------------------------------------------------------

$marshalRange = $d.Content.Find.Execute(
    "and then get sent through socket.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)

$marshalPara = $d.Paragraphs.Item(11)

$marshalPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item(12).Range.Text = "Unmarshalling received data is pretty much the same process in same order. First, it extracts length and msgType from the received message. Then based on the msgType, the system figures out which elements to extract from the message. For example, if the msgType is RegSuccess, the system knows that there will be a server_identifier and port marshalled in the message. Hence, it extracts these two elements based on their size. "

$d.Paragraphs.Item(12).Range.InsertParagraphAfter()
$d.Paragraphs.Item(13).Range.Text = "These extracted, or unmarshalled, data gets passed into FunctionData constructor in order to further unmarshall detailed argTypes and args. These data is either saved in local database or gets used in different purposes."

# ---------------------------------------------------------------------------
# 3. "Structure of binder database" paragraph: rename server_struct ->
#    ServerData and extend the description with the FunctionData class.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "server_struct. The structure contains basic server information like hostname and port number. In addition, it has another vector of function skeletons that the server provides. With this implementation, location request from clients and register request from servers can be effectively handled. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ServerData. The structure contains basic server information like hostname and port number. In addition, it has another vector of custom class called FunctionData. FunctionData class is basically the same manner as ServerData, which holds information regarding on remote functions. As rpcRegister message comes in from the server, the system checks whether a FunctionData constructed based on the sent name and argTypes exist in the database. If there is an existence, then the binder sends back a RegSuccessMessage back to the server with a reason code indicating that there is a duplicated function definition in the binder. Otherwise, it simply adds the new FunctionData to the corresponding vector.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "Handling of function overloading" paragraph: "As a rpcRegister()" ->
#    "As rpcRegister()" (confine the edit to the run ending "...As a " so the
#    proofErr markers around the following "rpcRegister(" run are untouched)
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "corresponding function skeleton. As a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "corresponding function skeleton. As ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "Managing round-robin scheduling" paragraph: add "global" and append a
#    new closing sentence about the database doubling as a queue.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "uses a single pointer and the pointer iterates through the existing list ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "uses a single global pointer and the pointer iterates through the existing list. Since our database is implemented in a server manner (i.e. elements within the vector is in form of ServerData), the database itself can be used as a queue.",
    2) | Out-Null
